$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 3
$ws.Range("B2").Value = 4

$ws.Range("G5").Select()
